$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Swap the F:V content of row 95 and row 96 (columns A-E, the index/
#    pais/torneio/temporada/data_partida stay put - only the match details
#    change place).
# ---------------------------------------------------------------------------
$row95 = @("Grecia", 2, "Liberia", 3, 2.71, "19/10/2023 18:43", 2.39, "22/10/2023 23:51", 3.38, "19/10/2023 18:43", 3.6, "22/10/2023 23:50", 2.57, "19/10/2023 18:43", 2.87, "22/10/2023 23:51", "https://www.betexplorer.com/football/costa-rica/primera-division/grecia-liberia/0lhAqXdC/")
$row96 = @("AD Santos", 2, "Zeledon", 0, 1.88, "19/10/2023 18:43", 1.97, "22/10/2023 23:50", 3.6, "19/10/2023 18:43", 3.52, "22/10/2023 23:50", 4.06, "19/10/2023 18:43", 3.95, "22/10/2023 23:50", "https://www.betexplorer.com/football/costa-rica/primera-division/santos-de-guapiles-zeledon/Sb0cnZCg/")

$cols = @("F", "G", "H", "I", "J", "K", "L", "M", "N", "O", "P", "Q", "R", "S", "T", "U", "V")

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "95").Value = $row95[$i]
    $ws.Range($cols[$i] + "96").Value = $row96[$i]
}

# ---------------------------------------------------------------------------
# 2) Append a new row 122 with a fresh match, cloning the cell formatting
#    used in row 121 (bold/bordered/centered index column A, custom
#    date-time format column E) via PasteSpecial so the same style indexes
#    get reused instead of creating brand new styles.
# ---------------------------------------------------------------------------
$ws.Range("A121").Copy()
$ws.Range("A122").PasteSpecial(-4122)
$ws.Range("E121").Copy()
$ws.Range("E122").PasteSpecial(-4122)

$ws.Range("A122").Value = 121
$ws.Range("B122").Value = "costa-rica"
$ws.Range("C122").Value = "primera-division"
$ws.Range("D122").Value = "2023-2024"
$ws.Range("E122").Value = 45252.125
$ws.Range("F122").Value = "Sporting San Jose"
$ws.Range("G122").Value = 1
$ws.Range("H122").Value = "Cartagines"
$ws.Range("I122").Value = 1
$ws.Range("J122").Value = 2.41
$ws.Range("K122").Value = "15/11/2023 05:12"
$ws.Range("L122").Value = 2.6
$ws.Range("M122").Value = "21/11/2023 22:48"
$ws.Range("N122").Value = 3.42
$ws.Range("O122").Value = "15/11/2023 05:12"
$ws.Range("P122").Value = 3.39
$ws.Range("Q122").Value = "21/11/2023 22:48"
$ws.Range("R122").Value = 2.88
$ws.Range("S122").Value = "15/11/2023 05:12"
$ws.Range("T122").Value = 2.49
$ws.Range("U122").Value = "21/11/2023 22:48"
$ws.Range("V122").Value = "https://www.betexplorer.com/football/costa-rica/primera-division/sporting-san-jose-cartagines/jZlw4bvd/"

Write-Output "done"
